# Auto-generated edits applying the Behemoth_Profits market-data refresh
# (scheduled runner price/profit update) described by the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 5000
$ws.Range("J13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("N13").Value = -5338

$ws.Range("H39").Value = 198.46666
$ws.Range("I39").Value = 195.75
$ws.Range("J39").Value = 209.33333
$ws.Range("K39").Value = 587.25
$ws.Range("L39").Value = 627.99999
$ws.Range("M39").Value = -291.25
$ws.Range("N39").Value = -1219.99999

$ws.Range("H40").Value = 3448.8
$ws.Range("I40").Value = 2824.75
$ws.Range("K40").Value = 2824.75
$ws.Range("M40").Value = -2649.75

$ws.Range("H69").Value = 31873.125
$ws.Range("J69").Value = 42502.5
$ws.Range("L69").Value = 127507.5
$ws.Range("N69").Value = -129255.5

$ws.Range("H72").Value = 31873.125
$ws.Range("J72").Value = 42502.5
$ws.Range("L72").Value = 382522.5
$ws.Range("N72").Value = -391258.5

$ws.Range("H105").Value = 51492.2
$ws.Range("J105").Value = 51492.2
$ws.Range("L105").Value = 51492.2
$ws.Range("N105").Value = -58480.2

$ws.Range("H137").Value = 3212.054
$ws.Range("I137").Value = 2277.4482
$ws.Range("J137").Value = 6600
$ws.Range("K137").Value = 6832.344599999999
$ws.Range("L137").Value = 19800
$ws.Range("M137").Value = -4282.344599999999
$ws.Range("N137").Value = -24900

$ws.Range("H138").Value = 2655.889
$ws.Range("I138").Value = 1439.4286
$ws.Range("J138").Value = 3205.258
$ws.Range("K138").Value = 4318.2858
$ws.Range("L138").Value = 9615.773999999999
$ws.Range("M138").Value = 821.7142000000003
$ws.Range("N138").Value = -19895.774

$ws.Range("H141").Value = 2438.318
$ws.Range("I141").Value = 2438.7058
$ws.Range("K141").Value = 7316.117400000001
$ws.Range("M141").Value = -2136.117400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5815896.5
$ws.Range("I32").Value = 5883848
$ws.Range("K32").Value = 5883848
$ws.Range("M32").Value = -5883561

$ws.Range("H94").Value = 46954.5
$ws.Range("J94").Value = 46954.5
$ws.Range("L94").Value = 46954.5
$ws.Range("N94").Value = -48756.5

$ws.Range("H106").Value = 50786.668
$ws.Range("J106").Value = 50786.668
$ws.Range("L106").Value = 50786.668
$ws.Range("N106").Value = -53310.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 66666
$ws.Range("J57").Value = 66666
$ws.Range("L57").Value = 66666
$ws.Range("N57").Value = -68106

$ws.Range("H135").Value = 32999.668
$ws.Range("J135").Value = 32999.668
$ws.Range("L135").Value = 32999.668
$ws.Range("N135").Value = -43139.668

$ws.Range("H136").Value = 66666
$ws.Range("J136").Value = 66666
$ws.Range("L136").Value = 66666
$ws.Range("N136").Value = -76866

$ws.Range("H137").Value = 70779
$ws.Range("J137").Value = 70779
$ws.Range("L137").Value = 70779
$ws.Range("N137").Value = -80979

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 514530.3
$ws.Range("I31").Value = 9488.869000000001
$ws.Range("J31").Value = 1019571.75
$ws.Range("K31").Value = 9488.869000000001
$ws.Range("L31").Value = 1019571.75
$ws.Range("M31").Value = -9193.869000000001
$ws.Range("N31").Value = -1020161.75

$ws.Range("H34").Value = 514530.3
$ws.Range("I34").Value = 9488.869000000001
$ws.Range("J34").Value = 1019571.75
$ws.Range("K34").Value = 9488.869000000001
$ws.Range("L34").Value = 1019571.75
$ws.Range("M34").Value = -9286.869000000001
$ws.Range("N34").Value = -1019975.75

$ws.Range("H47").Value = 57535.5
$ws.Range("I47").Value = 25000
$ws.Range("K47").Value = 25000
$ws.Range("M47").Value = -24434

$ws.Range("H53").Value = 41007.832
$ws.Range("J53").Value = 41007.832
$ws.Range("L53").Value = 41007.832
$ws.Range("N53").Value = -42221.832

$ws.Range("H58").Value = 1555.0952
$ws.Range("I58").Value = 1297.3572
$ws.Range("J58").Value = 2070.5715
$ws.Range("K58").Value = 1297.3572
$ws.Range("L58").Value = 2070.5715
$ws.Range("M58").Value = -1094.3572
$ws.Range("N58").Value = -2476.5715

$ws.Range("H95").Value = 54944.5
$ws.Range("J95").Value = 54944.5
$ws.Range("L95").Value = 54944.5
$ws.Range("N95").Value = -60436.5

$ws.Range("H99").Value = 6860.579
$ws.Range("I99").Value = 6785
$ws.Range("J99").Value = 7263.6665
$ws.Range("K99").Value = 6785
$ws.Range("L99").Value = 7263.6665
$ws.Range("M99").Value = -5287
$ws.Range("N99").Value = -10259.6665

$ws.Range("H107").Value = 953.6667
$ws.Range("I107").Value = 439.46155
$ws.Range("J107").Value = 1789.25
$ws.Range("K107").Value = 439.46155
$ws.Range("L107").Value = 1789.25
$ws.Range("M107").Value = 1480.53845
$ws.Range("N107").Value = -5629.25

$ws.Range("H126").Value = 6860.579
$ws.Range("I126").Value = 6785
$ws.Range("J126").Value = 7263.6665
$ws.Range("K126").Value = 20355
$ws.Range("L126").Value = 21790.9995
$ws.Range("M126").Value = -17885
$ws.Range("N126").Value = -26730.9995

$ws.Range("H134").Value = 313857.47
$ws.Range("I134").Value = 455488.38
$ws.Range("J134").Value = 2269.5
$ws.Range("K134").Value = 1366465.14
$ws.Range("L134").Value = 6808.5
$ws.Range("M134").Value = -1363930.14
$ws.Range("N134").Value = -11878.5

$ws.Range("H136").Value = 1555.0952
$ws.Range("I136").Value = 1297.3572
$ws.Range("J136").Value = 2070.5715
$ws.Range("K136").Value = 3892.0716
$ws.Range("L136").Value = 6211.7145
$ws.Range("M136").Value = -1342.0716
$ws.Range("N136").Value = -11311.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2553.52
$ws.Range("I11").Value = 2534.9167
$ws.Range("J11").Value = 3000
$ws.Range("K11").Value = 7604.750100000001
$ws.Range("L11").Value = 9000
$ws.Range("M11").Value = -7464.750100000001
$ws.Range("N11").Value = -9280

$ws.Range("H132").Value = 1907.5264
$ws.Range("I132").Value = 2006.1
$ws.Range("K132").Value = 18054.9
$ws.Range("M132").Value = -15524.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 18000
$ws.Range("J26").Value = 18000
$ws.Range("L26").Value = 18000
$ws.Range("N26").Value = -18560

$ws.Range("H50").Value = 18000
$ws.Range("J50").Value = 18000
$ws.Range("L50").Value = 18000
$ws.Range("N50").Value = -18996

$ws.Range("H53").Value = 25000
$ws.Range("J53").Value = 25000
$ws.Range("L53").Value = 25000
$ws.Range("N53").Value = -26262

$ws.Range("H58").Value = 25000
$ws.Range("J58").Value = 25000
$ws.Range("L58").Value = 25000
$ws.Range("N58").Value = -25554

$ws.Range("H95").Value = 125039670
$ws.Range("J95").Value = 125039670
$ws.Range("L95").Value = 125039670
$ws.Range("N95").Value = -125045162

$ws.Range("H102").Value = 3075.838
$ws.Range("I102").Value = 2466.8333
$ws.Range("J102").Value = 25000
$ws.Range("K102").Value = 2466.8333
$ws.Range("L102").Value = 25000
$ws.Range("M102").Value = -844.8332999999998
$ws.Range("N102").Value = -28244

$ws.Range("H132").Value = 34488490
$ws.Range("I132").Value = 55557084
$ws.Range("J132").Value = 12609
$ws.Range("K132").Value = 166671252
$ws.Range("L132").Value = 37827
$ws.Range("M132").Value = -166668722
$ws.Range("N132").Value = -42887

$ws.Range("H140").Value = 78779.39999999999
$ws.Range("J140").Value = 78779.39999999999
$ws.Range("L140").Value = 78779.39999999999
$ws.Range("N140").Value = -89139.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3471.923
$ws.Range("I40").Value = 3471.923
$ws.Range("K40").Value = 3471.923
$ws.Range("M40").Value = -3335.923

$ws.Range("H100").Value = 4358
$ws.Range("I100").Value = 3042.3333
$ws.Range("J100").Value = 6331.5
$ws.Range("K100").Value = 3042.3333
$ws.Range("L100").Value = 6331.5
$ws.Range("M100").Value = -2501.3333
$ws.Range("N100").Value = -7413.5

$ws.Range("H103").Value = 57401.332
$ws.Range("J103").Value = 57401.332
$ws.Range("L103").Value = 57401.332
$ws.Range("N103").Value = -59745.332

$ws.Range("H132").Value = 19757.188
$ws.Range("I132").Value = 4808.452
$ws.Range("K132").Value = 14425.356
$ws.Range("M132").Value = -11895.356

$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2405.5
$ws.Range("I122").Value = 1817.8
$ws.Range("K122").Value = 5453.4
$ws.Range("M122").Value = -3003.4

$ws.Range("H132").Value = 2330.1538
$ws.Range("I132").Value = 2327.36
$ws.Range("K132").Value = 6982.08
$ws.Range("M132").Value = -4452.08

$ws.Range("H136").Value = 3787.5
$ws.Range("I136").Value = 718.4545000000001
$ws.Range("K136").Value = 2155.3635
$ws.Range("M136").Value = 394.6364999999996
